$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.688.07"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.887.40"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -1.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4847"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07325"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9182"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07676"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "1.916.33"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.459"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.587"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008787"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "27.722.81"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "2.144.68"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.910"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.106"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.890"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08930"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.154"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.218"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.621"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02036"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -6.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.090"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.00%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05248"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5442"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.973"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.933"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1516"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4769"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.628"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06053"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.80%  "
